$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 2, so the old row2 content shifts to row3
$ws.Rows.Item(2).Insert()

# Row 3: remainder of the text + the data that was already there
$ws.Range("A3").Value = "expiration time"

# Row 2: split "Access Token expiration time" into A2 = "Access Token " (highlighted)
$ws.Range("A2").Style = "Good"
$ws.Range("A2").Value = "Access Token "

# Update selection to match target
$ws.Range("C8").Select()
